# "Add files via upload" — the Bookings sheet was re-uploaded as a fresh
# template: the two sample booking rows are gone, a new "Status" column
# was appended after "Confirmation Number", and a handful of blank
# (pre-formatted) rows were left under the header for future entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bookings")

# Drop the two sample data rows (rows 2 and 3), keeping just the header.
$ws.Rows("2:3").Delete()

# Append a new "Status" column (J) after "Confirmation Number" (I).
$ws.Columns("J:J").Insert()
$ws.Range("J1").Value = "Status"

# Leave the cursor/selection parked the way the re-uploaded sheet has it:
# a block of blank rows below the header, ready for new bookings.
$ws.Range("A2:K6").Select()
